$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data: Date 45979 (2025-11-18), Error Count 2
# Copy the formatting from the last existing row (A14) so the new date
# cell (A15) picks up the same date-number-format style without minting
# a brand new style record.
$ws.Range("A14").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A15").Value = 45979
$ws.Range("B15").Value = 2

# Update selection to match the new last row, like Excel would leave it
# after the data entry.
$ws.Range("A15:B15").Select()
